$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 77
$ws1.Range("F4").Value = 2094
$ws1.Range("F5").Value = 371
$ws1.Range("F6").Value = 652
$ws1.Range("F8").Value = 2084
$ws1.Range("F9").Value = 10774
$ws1.Range("F15").Value = 8999
$ws1.Range("F16").Value = 1117
$ws1.Range("F17").Value = 730
$ws1.Range("F18").Value = 5282
$ws1.Range("F20").Value = 3354

# Sheet "全部类型" (sheet4) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 77
$ws4.Range("F4").Value = 2094
$ws4.Range("F5").Value = 371
$ws4.Range("F6").Value = 652
$ws4.Range("F9").Value = 2084
$ws4.Range("F12").Value = 10774
$ws4.Range("F18").Value = 8999
$ws4.Range("F19").Value = 1117
$ws4.Range("F20").Value = 730
$ws4.Range("F21").Value = 5282
$ws4.Range("F23").Value = 3354
